$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(14.83701644944724,0,11.24699344450864,16.89970777792428,28.18381175807251,26.55036001338741,13.52405017548072,24.31714571321515,11.62502896645641,0,9.827231030235167,15.09242133468871,18.78675607785098,20.32790489034049),
    @(14.57636004658343,0,11.27502430767402,16.90984672940891,28.21924535121965,26.39446825181204,13.54779962693691,24.44654682276348,11.62090311551327,0,9.711904924517501,14.96268961216196,18.8035701488626,20.33188376275236),
    @(14.41532735935859,0,11.29323882438274,16.91717580308636,28.24786014031482,26.30720427913595,13.56474931343558,24.5300356459713,11.61882375850214,0,9.64119651565156,14.88377829207912,18.81579662716122,20.33910700399061),
    @(14.34953458766999,0,11.30091435555018,16.9204411250671,28.26124345252709,26.27380435563896,13.57225152717519,24.5650756488761,11.61809140463989,0,9.612434809273605,14.85183438432961,18.82125867469346,20.34325301307465),
    @(14.33860175122498,0,11.30220416917126,16.92100019268637,28.26356970587648,26.26838977757827,13.57353319671965,24.57095555765692,11.61797676841867,0,9.607662906483901,14.84654374233936,18.82219464931301,20.34401408207221),
    @(14.41444064178949,0,11.29334131416053,16.91721871056715,28.24803366101816,26.30674504364439,13.56484808184454,24.53050408371347,11.61881341496152,0,9.640808377069106,14.88334658864077,18.81586834650988,20.33915804961629),
    @(14.74739352855726,0,11.25645060573775,16.9029752906642,28.19460476203997,26.49487606194614,13.53174751567733,24.36092719508436,11.62351267455695,0,9.787461449771181,15.04755095424521,18.79215944463661,20.32828514280896),
    @(15.38914223499363,0,11.19204083461184,16.88375209207289,28.14431360429498,26.92915430502664,13.48562929870548,24.06028504128787,11.63629249060334,0,10.07468478935345,15.37431242604488,18.76070762220257,20.34484550505893),
    @(15.84953908691048,0,11.14951435875633,16.87487125669278,28.14062577505364,27.28549445709483,13.46320807789023,23.85867004487232,11.64780561353245,0,10.28392864770607,15.61573347652439,18.74669549272877,20.38000453719854),
    @(16.05572574915369,0,11.13120061872054,16.87195513062583,28.1461648399225,27.45508368295715,13.4554960045371,23.77109569900785,11.65349369796175,0,10.37840947079751,15.7255247053511,18.74227951410612,20.40095564961137),
    @(16.13326937309173,0,11.12441339892636,16.87101123418201,28.14929819817226,27.52032408535676,13.45293303668249,23.73852627074617,11.65571148803101,0,10.41405912796539,15.76706748608507,18.74088740428871,20.40959782294935),
    @(16.11659378615917,0,11.12586858390928,16.87120740604997,28.14857734367788,27.5062289212093,13.45346912656541,23.74551434127721,11.65523102382125,0,10.40638745581034,15.75812240708789,18.74117478519435,20.40770514564675),
    @(16.06211640107287,0,11.13063927125791,16.87187426857897,28.14640187742975,27.46043092432681,13.45527798622608,23.76840432031719,11.65367488315524,0,10.38134510923275,15.72894327176243,18.74215937881097,20.40165247586441),
    @(16.02867584494376,0,11.13358068558226,16.87230359031749,28.14520416684332,27.43250954252354,13.45643250252495,23.78250225735511,11.65272998619873,0,10.36598847383003,15.71106511118501,18.74279890720391,20.39803715844049),
    @(15.83599305668185,0,11.15073191761246,16.87508434197926,28.14040885730365,27.27455733609607,13.46376211405011,23.86447635911347,11.64744287018915,0,10.27773760828837,15.60855529971079,18.7470233675575,20.37873472410897),
    @(15.71690685693691,0,11.16151749534973,16.87707721152724,28.13931435383089,27.17953747165972,13.46889553989937,23.91582376330679,11.64431414345808,0,10.22339866656318,15.5456404617744,18.75011545604647,20.36816004929661),
    @(15.64810951235427,0,11.16781822300087,16.87832931014889,28.1393642588742,27.12559291071812,13.47208230514309,23.94574737784749,11.6425570645383,0,10.19207952249705,15.5094521276715,18.75207834998746,20.36254466951582),
    @(15.62476613844289,0,11.16996824396343,16.87877146658689,28.13949789423285,27.10745154366091,13.47320151142727,23.95594604794697,11.64196947359261,0,10.18146509205751,15.49720000296152,18.75277466791431,20.36072370943493),
    @(15.72961555343617,0,11.16035930183163,16.87685412062721,28.13936055791556,27.18957958645476,13.46832484536879,23.91031740431891,11.64464281205865,0,10.22919004460104,15.5523382013692,18.74976722301884,20.36923744918538),
    @(16.0781327630437,0,11.12923399811878,16.87167405215075,28.14701277170685,27.47385566455841,13.45473698256706,23.76166490204073,11.65413023474561,0,10.38870434198612,15.73751500320558,18.74186258916993,20.40341110331261),
    @(16.30276304308056,0,11.10975302072276,16.86922279824457,28.15805009923668,27.66556810790421,13.44793967472874,23.66796768320614,11.6607024359804,0,10.49219713305588,15.8583348994794,18.7383285268488,20.42987235900175),
    @(16.18318281519578,0,11.12007177152686,16.87044600669163,28.15160780697861,27.56272483855698,13.45137703797448,23.71766024812073,11.65716105097912,0,10.43703918683482,15.79387890594026,18.74006588331957,20.41537352149748),
    @(15.72387098410747,0,11.160882609742,16.87695464860918,28.13933755341641,27.1850374126701,13.46858212280407,23.91280557275957,11.64449409112479,0,10.2265720064895,15.54931020644465,18.74992408214586,20.36874891059811),
    @(15.21718137393312,0,11.20862034775239,16.88802699188281,28.15207806801944,26.80492557496258,13.49609265737704,24.13822042872009,11.63245862912303,0,10.07468478935345,15.28556523368316,18.76761337479058,20.33631765764735)
)

for ($i = 0; $i -lt 24; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt 14; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}
